$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold values (row 2 / row 3, columns B & C)
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 9.5

# Move the active selection from C3 to B3, matching the saved view state
$ws.Range("B3").Select()
